# Automatic update of files.
#
# The source export re-ran and reshuffled the occurrence rows (the report's
# "Id" values moved between rows, coordinates got rounded to whole metres,
# and the now-unused Starttid/Sluttid ("00:00") columns were dropped), plus
# a couple of public-comment text tweaks. Apply the new values cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = 111639170
$ws.Range("Q2").Value = 548231
$ws.Range("R2").Value = 6926520
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = 111639167
$ws.Range("Q3").Value = 547815
$ws.Range("R3").Value = 6926124
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
$ws.Range("AC3").Value = "1 planta"

# --- Row 4 -----------------------------------------------------------------
$ws.Range("A4").Value = 111639175
$ws.Range("B4").Value = 89686
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 658
$ws.Range("F4").Value = "Rosenticka"
$ws.Range("G4").Value = "Rhodofomes roseus"
$ws.Range("H4").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q4").Value = 547828
$ws.Range("R4").Value = 6926125
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = 111639168
$ws.Range("B5").Value = 89686
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 658
$ws.Range("F5").Value = "Rosenticka"
$ws.Range("G5").Value = "Rhodofomes roseus"
$ws.Range("H5").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q5").Value = 548104
$ws.Range("R5").Value = 6926478
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
$ws.Range("AC5").ClearContents()

# --- Row 6 (same occurrence, only coordinates rounded + time cols gone) --
$ws.Range("Q6").Value = 548225
$ws.Range("R6").Value = 6926513
$ws.Range("Z6").ClearContents()
$ws.Range("AB6").ClearContents()

# --- Row 7 -----------------------------------------------------------------
$ws.Range("A7").Value = 111639174
$ws.Range("B7").Value = 96348
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("Q7").Value = 547804
$ws.Range("R7").Value = 6926147
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").Value = "ca 6 plantor"

# --- Row 8 ---------------------------------------------------------------
$ws.Range("A8").Value = 111639173
$ws.Range("Q8").Value = 547838
$ws.Range("R8").Value = 6926229
$ws.Range("Z8").ClearContents()
$ws.Range("AB8").ClearContents()
$ws.Range("AC8").Value = "ca 15 plantor"

# --- Row 9 -----------------------------------------------------------------
$ws.Range("A9").Value = 111639172
$ws.Range("B9").Value = 96348
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = "Knärot"
$ws.Range("G9").Value = "Goodyera repens"
$ws.Range("H9").Value = "(L.) R. Br."
$ws.Range("Q9").Value = 548221
$ws.Range("R9").Value = 6926512
$ws.Range("Z9").ClearContents()
$ws.Range("AB9").ClearContents()
